$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy style (border/bold/alignment) from row 6 down into the new rows 7-11 first,
# so the values we set afterwards keep their own cell refs/values untouched. ---
$ws.Range("A6").Copy($ws.Range("A7:A11"))

# --- Set B7:B11 (new song titles) first so the new unique strings are appended
# to the shared-string table in this order (indices 4..8). ---
$ws.Range("B7").Value = 'Hold On, We''re Going Home'
$ws.Range("B8").Value = 'I Took A Pill In Ibiza - Seeb Remix'
$ws.Range("B9").Value = 'Love Yourself'
$ws.Range("B10").Value = 'Side To Side'
$ws.Range("B11").Value = 'Truffle Butter'

# --- Now overwrite C2:C5 with the re-serialised (key-reordered) JSON text.
# These become new unique strings appended right after the titles (indices 9..12),
# and the old JSON strings they replace get garbage-collected on save. ---
$ws.Range("C2").Value = '[{''key'': 7, ''mode'': 1, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/2IY7eOUDjw2ArKYxKa2jXc'', ''valence'': 0.477, ''energy'': 0.592, ''liveness'': 0.136, ''speechiness'': 0.281, ''loudness'': -7.033, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/2IY7eOUDjw2ArKYxKa2jXc'', ''danceability'': 0.682, ''acousticness'': 0.169, ''tempo'': 186.041, ''duration_ms'': 230467, ''id'': ''2IY7eOUDjw2ArKYxKa2jXc'', ''instrumentalness'': 5.49e-06, ''time_signature'': 4, ''uri'': ''spotify:track:2IY7eOUDjw2ArKYxKa2jXc''}]'
$ws.Range("C3").Value = '[{''key'': 1, ''mode'': 1, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/12VWzyPDBCc8fqeWCAfNwR'', ''valence'': 0.378, ''energy'': 0.617, ''liveness'': 0.351, ''speechiness'': 0.0522, ''loudness'': -5.871, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/12VWzyPDBCc8fqeWCAfNwR'', ''danceability'': 0.785, ''acousticness'': 0.00902, ''tempo'': 103.981, ''duration_ms'': 173987, ''id'': ''12VWzyPDBCc8fqeWCAfNwR'', ''instrumentalness'': 0.00246, ''time_signature'': 4, ''uri'': ''spotify:track:12VWzyPDBCc8fqeWCAfNwR''}]'
$ws.Range("C4").Value = '[{''key'': 7, ''mode'': 1, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/7fJtPlEZKxu6gvkfBFc5tW'', ''valence'': 0.391, ''energy'': 0.65, ''liveness'': 0.102, ''speechiness'': 0.118, ''loudness'': -7.79, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/7fJtPlEZKxu6gvkfBFc5tW'', ''danceability'': 0.804, ''acousticness'': 0.0606, ''tempo'': 117.984, ''duration_ms'': 263373, ''id'': ''7fJtPlEZKxu6gvkfBFc5tW'', ''instrumentalness'': 7.05e-05, ''time_signature'': 4, ''uri'': ''spotify:track:7fJtPlEZKxu6gvkfBFc5tW''}]'
$ws.Range("C5").Value = '[{''key'': 3, ''mode'': 1, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/7LVHVU3tWfcxj5aiPFEW4Q'', ''valence'': 0.122, ''energy'': 0.418, ''liveness'': 0.113, ''speechiness'': 0.0338, ''loudness'': -8.74, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/7LVHVU3tWfcxj5aiPFEW4Q'', ''danceability'': 0.209, ''acousticness'': 0.163, ''tempo'': 138.265, ''duration_ms'': 295533, ''id'': ''7LVHVU3tWfcxj5aiPFEW4Q'', ''instrumentalness'': 0.00195, ''time_signature'': 4, ''uri'': ''spotify:track:7LVHVU3tWfcxj5aiPFEW4Q''}]'

# --- C6 duplicates the (now updated) C4 value. ---
$ws.Range("C6").Value = '[{''key'': 7, ''mode'': 1, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/7fJtPlEZKxu6gvkfBFc5tW'', ''valence'': 0.391, ''energy'': 0.65, ''liveness'': 0.102, ''speechiness'': 0.118, ''loudness'': -7.79, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/7fJtPlEZKxu6gvkfBFc5tW'', ''danceability'': 0.804, ''acousticness'': 0.0606, ''tempo'': 117.984, ''duration_ms'': 263373, ''id'': ''7fJtPlEZKxu6gvkfBFc5tW'', ''instrumentalness'': 7.05e-05, ''time_signature'': 4, ''uri'': ''spotify:track:7fJtPlEZKxu6gvkfBFc5tW''}]'

# --- Fill in the rest of the new rows 7-11. ---
$ws.Range("A7").Value = 5
$ws.Range("C7").Value = '[{''key'': 6, ''mode'': 0, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/14Rcq31SafFBHNEwXrtR2B'', ''valence'': 0.287, ''energy'': 0.407, ''liveness'': 0.0705, ''speechiness'': 0.0907, ''loudness'': -7.442, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/14Rcq31SafFBHNEwXrtR2B'', ''danceability'': 0.772, ''acousticness'': 0.00349, ''tempo'': 100.006, ''duration_ms'': 227880, ''id'': ''14Rcq31SafFBHNEwXrtR2B'', ''instrumentalness'': 2.4e-05, ''time_signature'': 4, ''uri'': ''spotify:track:14Rcq31SafFBHNEwXrtR2B''}]'
$ws.Range("A8").Value = 6
$ws.Range("C8").Value = '[{''key'': 7, ''mode'': 0, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/1MtUq6Wp1eQ8PC6BbPCj8P'', ''valence'': 0.697, ''energy'': 0.725, ''liveness'': 0.0864, ''speechiness'': 0.12, ''loudness'': -6.588, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/1MtUq6Wp1eQ8PC6BbPCj8P'', ''danceability'': 0.665, ''acousticness'': 0.034, ''tempo'': 101.96, ''duration_ms'': 197933, ''id'': ''1MtUq6Wp1eQ8PC6BbPCj8P'', ''instrumentalness'': 9.17e-06, ''time_signature'': 4, ''uri'': ''spotify:track:1MtUq6Wp1eQ8PC6BbPCj8P''}]'
$ws.Range("A9").Value = 7
$ws.Range("C9").Value = '[{''key'': 4, ''mode'': 1, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/50kpGaPAhYJ3sGmk6vplg0'', ''valence'': 0.559, ''energy'': 0.376, ''liveness'': 0.285, ''speechiness'': 0.453, ''loudness'': -9.954, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/50kpGaPAhYJ3sGmk6vplg0'', ''danceability'': 0.607, ''acousticness'': 0.856, ''tempo'': 102.541, ''duration_ms'': 233720, ''id'': ''50kpGaPAhYJ3sGmk6vplg0'', ''instrumentalness'': 0, ''time_signature'': 4, ''uri'': ''spotify:track:50kpGaPAhYJ3sGmk6vplg0''}]'
$ws.Range("A10").Value = 8
$ws.Range("C10").Value = '[{''key'': 6, ''mode'': 0, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/1pKeFVVUOPjFsOABub0OaV'', ''valence'': 0.624, ''energy'': 0.728, ''liveness'': 0.328, ''speechiness'': 0.232, ''loudness'': -5.914, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/1pKeFVVUOPjFsOABub0OaV'', ''danceability'': 0.65, ''acousticness'': 0.0373, ''tempo'': 159.167, ''duration_ms'': 226160, ''id'': ''1pKeFVVUOPjFsOABub0OaV'', ''instrumentalness'': 0, ''time_signature'': 4, ''uri'': ''spotify:track:1pKeFVVUOPjFsOABub0OaV''}]'
$ws.Range("A11").Value = 9
$ws.Range("C11").Value = '[{''key'': 10, ''mode'': 0, ''type'': ''audio_features'', ''track_href'': ''https://api.spotify.com/v1/tracks/3keUgTGEoZJt0QkzTB6kHg'', ''valence'': 0.511, ''energy'': 0.672, ''liveness'': 0.12, ''speechiness'': 0.0479, ''loudness'': -6.875, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/3keUgTGEoZJt0QkzTB6kHg'', ''danceability'': 0.896, ''acousticness'': 0.0525, ''tempo'': 105.098, ''duration_ms'': 219227, ''id'': ''3keUgTGEoZJt0QkzTB6kHg'', ''instrumentalness'': 3.14e-05, ''time_signature'': 4, ''uri'': ''spotify:track:3keUgTGEoZJt0QkzTB6kHg''}]'
